$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from H1 (bold, bordered, centered) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0/IF data for data rows 2-72 (index 0 corresponds to row 2)
$i0 = @(8,7,7,8,11,6,6,3,4,8,8,7,8,7,9,7,5,6,9,7,3,9,9,9,9,8,8,9,10,8,8,8,7,8,8,6,8,7,8,5,9,6,7,7,8,6,8,6,7,7,6,7,6,8,7,6,10,9,7,9,8,8,7,5,5,5,5,8,4,8,5)
$if = @(9,7,7,8,12,6,8,4,5,8,8,7,9,7,9,8,6,7,9,8,4,9,9,9,9,9,8,9,10,9,8,8,7,8,8,6,8,8,8,6,9,6,8,7,9,6,9,7,7,7,7,7,7,8,7,6,10,9,7,9,8,8,8,6,5,5,5,8,4,8,5)

for ($n = 0; $n -lt $i0.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $i0[$n]
    $ws.Cells.Item($row, 10).Value = $if[$n]
}
